$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.722.98'
$ws.Cells.Item(2, 5).Value = '  -0.42%  '
$ws.Cells.Item(3, 4).Value = '1.632.05'
$ws.Cells.Item(3, 5).Value = '  -0.21%  '
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '214.92'
$ws.Cells.Item(5, 5).Value = '  -0.08%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '0.504'
$ws.Cells.Item(6, 5).Value = '  -0.88%  '
$ws.Cells.Item(7, 5).Value = '  -0.11%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 5).Value = '  -1.13%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '19.65'
$ws.Cells.Item(10, 5).Value = '  -3.07%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.0786'
$ws.Cells.Item(11, 5).Value = '  +0.87%  '
$ws.Cells.Item(12, 5).Value = '  -0.20%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.633.66'
$ws.Cells.Item(13, 5).Value = '  -0.35%  '
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '1.856.27'
$ws.Cells.Item(14, 5).Value = '  -0.25%  '
$ws.Cells.Item(15, 5).Value = '  -0.39%  '
$ws.Cells.Item(16, 5).Value = '  -0.39%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '62.64'
$ws.Cells.Item(17, 5).Value = '  -0.92%  '
$ws.Cells.Item(18, 4).Value = '25.745.14'
$ws.Cells.Item(18, 5).Value = '  -0.35%  '
$ws.Cells.Item(20, 5).Value = '  +1.62%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '193.79'
$ws.Cells.Item(21, 5).Value = '  +0.90%  '
$ws.Cells.Item(22, 5).Value = '  +0.10%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '6.26'
$ws.Cells.Item(23, 5).Value = '  +2.00%  '
$ws.Cells.Item(24, 5).Value = '  -0.14%  '
$ws.Cells.Item(25, 5).Value = '  +4.23%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '142.66'
$ws.Cells.Item(26, 5).Value = '  +2.66%  '
$ws.Cells.Item(27, 5).Value = '  -0.32%  '
$ws.Cells.Item(28, 5).Value = '  +0.62%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '15.51'
$ws.Cells.Item(29, 5).Value = '  -0.32%  '
$ws.Cells.Item(30, 5).Value = '  -0.19%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '0.0492'
$ws.Cells.Item(31, 5).Value = '  -0.64%  '
$ws.Cells.Item(32, 5).Value = '  +0.64%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '3.24'
$ws.Cells.Item(33, 5).Value = '  -0.57%  '
$ws.Cells.Item(34, 5).Value = '  +0.52%  '
$ws.Cells.Item(35, 5).Value = '  +0.03%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.899'
$ws.Cells.Item(36, 5).Value = '  -0.16%  '
$ws.Cells.Item(37, 4).Value = '1.125.47'
$ws.Cells.Item(37, 5).Value = '  -0.45%  '
$ws.Cells.Item(38, 5).Value = '  -1.56%  '
$ws.Cells.Item(39, 5).Value = '  -2.10%  '
$ws.Cells.Item(40, 5).Value = '  -1.00%  '
$ws.Cells.Item(42, 5).Value = '  +2.01%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '99.71'
$ws.Cells.Item(43, 5).Value = '  +0.69%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.803'
$ws.Cells.Item(44, 5).Value = '  +0.54%  '
$ws.Cells.Item(45, 4).Value = '1.766.84'
$ws.Cells.Item(45, 5).Value = '  -0.30%  '
$ws.Cells.Item(46, 5).Value = '  -1.23%  '
$ws.Cells.Item(47, 5).Value = '  -1.17%  '
$ws.Cells.Item(49, 5).Value = '  +0.13%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '7.59'
$ws.Cells.Item(50, 5).Value = '  -2.69%  '
$ws.Cells.Item(51, 2).Value = 'SynthetixNetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '2.34'
$ws.Cells.Item(51, 5).Value = '  +3.14%  '
